$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the text-value range (columns B:E, which hold coin name, link, price
# and volume as text) to Text format first, so numeric-looking strings such as
# "1.001" or "29.895.64" are written back as text rather than being parsed as
# numbers -- matching the original inlineStr cell type.
$rng = $ws.Range("B2:E51")
$rng.NumberFormat = "@"

$ws.Range("D2").Value = '29.895.64'
$ws.Range("E2").Value = '  +5.94%  '

$ws.Range("D3").Value = '1.877.88'

$ws.Range("E4").Value = '  +0.13%  '

$ws.Range("D5").Value = '248.20'
$ws.Range("E5").Value = '  +1.75%  '

$ws.Range("D7").Value = '0.4970'
$ws.Range("E7").Value = '  +1.07%  '

$ws.Range("D8").Value = '45.95'
$ws.Range("E8").Value = '  +9.21%  '

$ws.Range("D9").Value = '0.2840'
$ws.Range("E9").Value = '  +6.54%  '

$ws.Range("D10").Value = '0.06524'
$ws.Range("E10").Value = '  +4.39%  '

$ws.Range("D11").Value = '1.880.86'
$ws.Range("E11").Value = '  +5.62%  '

$ws.Range("D12").Value = '17.01'
$ws.Range("E12").Value = '  +3.10%  '

$ws.Range("D13").Value = '0.07228'
$ws.Range("E13").Value = '  +3.01%  '

$ws.Range("D14").Value = '0.6617'
$ws.Range("E14").Value = '  +5.62%  '

$ws.Range("D15").Value = '84.72'
$ws.Range("E15").Value = '  +6.16%  '

$ws.Range("D16").Value = '4.792'
$ws.Range("E16").Value = '  +3.33%  '

$ws.Range("D17").Value = '29.918.65'
$ws.Range("E17").Value = '  +6.12%  '

$ws.Range("D18").Value = '1.001'
$ws.Range("E18").Value = '  +0.03%  '

$ws.Range("D19").Value = '12.83'
$ws.Range("E19").Value = '  +6.42%  '

$ws.Range("D20").Value = '0.000007489'
$ws.Range("E20").Value = '  +3.59%  '

$ws.Range("D21").Value = '1.002'
$ws.Range("E21").Value = '  +0.21%  '

$ws.Range("D22").Value = '2.122.46'
$ws.Range("E22").Value = '  +5.75%  '

$ws.Range("D23").Value = '4.733'
$ws.Range("E23").Value = '  +3.77%  '

$ws.Range("D24").Value = '5.516'
$ws.Range("E24").Value = '  +5.48%  '

$ws.Range("D25").Value = '8.987'
$ws.Range("E25").Value = '  +2.88%  '

$ws.Range("D26").Value = '145.23'
$ws.Range("E26").Value = '  +3.13%  '

$ws.Range("D27").Value = '134.44'
$ws.Range("E27").Value = '  +23.33%  '

$ws.Range("D28").Value = '16.62'
$ws.Range("E28").Value = '  +5.46%  '

$ws.Range("D29").Value = '1.943'
$ws.Range("E29").Value = '  +4.77%  '

$ws.Range("D30").Value = '1.376'
$ws.Range("E30").Value = '  -0.65%  '

$ws.Range("D31").Value = '4.167'
$ws.Range("E31").Value = '  -0.06%  '

$ws.Range("D32").Value = '0.08586'
$ws.Range("E32").Value = '  +4.17%  '

$ws.Range("D33").Value = '3.852'
$ws.Range("E33").Value = '  +2.47%  '

$ws.Range("D34").Value = '0.05084'
$ws.Range("E34").Value = '  +4.02%  '

$ws.Range("D35").Value = '1.124'
$ws.Range("E35").Value = '  +5.01%  '

$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").Value = '0.6823'
$ws.Range("E36").Value = '  +4.98%  '

$ws.Range("B37").Value = 'HuobiToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D37").Value = '2.703'
$ws.Range("E37").Value = '  +3.42%  '

$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D38").Value = '2.289'
$ws.Range("E38").Value = '  +11.84%  '

$ws.Range("B39").Value = 'MXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D39").Value = '2.743'
$ws.Range("E39").Value = '  +5.77%  '

$ws.Range("B40").Value = 'TrustWalletToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D40").Value = '0.9567'
$ws.Range("E40").Value = '  +1.12%  '

$ws.Range("B41").Value = 'VeChain'
$ws.Range("C41").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D41").Value = '0.01625'
$ws.Range("E41").Value = '  +4.97%  '

$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").Value = '6.091'
$ws.Range("E42").Value = '  +3.18%  '

$ws.Range("B43").Value = 'Quant'
$ws.Range("C43").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D43").Value = '104.08'
$ws.Range("E43").Value = '  +4.88%  '

$ws.Range("B44").Value = 'PaxDollar'
$ws.Range("C44").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D44").Value = '1.001'
$ws.Range("E44").Value = '  +0.20%  '

$ws.Range("B45").Value = 'TheSandbox'
$ws.Range("C45").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D45").Value = '0.4195'
$ws.Range("E45").Value = '  +5.59%  '

$ws.Range("B46").Value = 'Aptos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D46").Value = '7.397'
$ws.Range("E46").Value = '  +3.23%  '

$ws.Range("B47").Value = 'Algorand'
$ws.Range("C47").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D47").Value = '0.1247'
$ws.Range("E47").Value = '  +3.61%  '

$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D48").Value = '0.05631'
$ws.Range("E48").Value = '  +3.71%  '

$ws.Range("B49").Value = 'Elrond'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D49").Value = '32.22'
$ws.Range("E49").Value = '  +5.25%  '

$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").Value = '8.230'
$ws.Range("E50").Value = '  +2.73%  '

$ws.Range("B51").Value = 'Decentraland'
$ws.Range("C51").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D51").Value = '0.3703'
$ws.Range("E51").Value = '  +6.86%  '

# Restore the default (Normal) style so no stray number-format styling remains
# on the cells we touched (keeps cell appearance identical to the original).
$rng.Style = "Normal"